# Generate Report for Archive
# Update the localization "Status" column from "Ready for handoff" to
# "In Translation" for the files that are now in translation
# (1fddbc87-7f0c-48a3-9e23-ab9824b42f19.md and
#  53857106-52f1-4d4d-bd0f-a9cfaba8a83b.md) on both the "zh-cn" and
# "de-de" language sheets.

$wb = $excel.ActiveWorkbook

$langSheets = @("zh-cn", "de-de")

foreach ($sheetName in $langSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    # Row 3 -> 1fddbc87-7f0c-48a3-9e23-ab9824b42f19.md
    # Row 4 -> 53857106-52f1-4d4d-bd0f-a9cfaba8a83b.md
    $ws.Range("C3").Value = "In Translation"
    $ws.Range("C4").Value = "In Translation"
}
